# Update localization status report to reflect "Ready for handoff" state,
# refreshing the handoff timestamps accordingly.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns (zh-cn / de-de) and Latest Handoff Date
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-28-19 18:28:53"

# zh-cn sheet: Status column and Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-19 18:28:50"

# de-de sheet: Status column and Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-19 18:28:53"
